# The deck's slide master theme ("Integral") and notes master theme
# ("Office Theme") are swapped: ppt/theme/theme1.xml should end up
# holding the "Office Theme" palette and ppt/theme/theme2.xml should
# end up holding the "Integral" palette.
#
# The PowerPoint object model reaches the deck's theme colour scheme
# through Slide.ThemeColorScheme, which maps 1:1 onto the 12 theme
# colour slots (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) of
# ppt/theme/theme1.xml.

function HexToRgbInt($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$officeTheme = @(
    "000000",  # 1  dk1
    "FFFFFF",  # 2  lt1
    "44546A",  # 3  dk2
    "E7E6E6",  # 4  lt2
    "5B9BD5",  # 5  accent1
    "ED7D31",  # 6  accent2
    "A5A5A5",  # 7  accent3
    "FFC000",  # 8  accent4
    "4472C4",  # 9  accent5
    "70AD47",  # 10 accent6
    "0563C1",  # 11 hlink
    "954F72"   # 12 folHlink
)

for ($i = 1; $i -le $officeTheme.Count; $i++) {
    $tcs.Colors($i).RGB = HexToRgbInt($officeTheme[$i - 1])
}
